$wb = $excel.ActiveWorkbook

# Sheet 1: question_answers - update column B (inline string answer values)
$ws1 = $wb.Worksheets.Item("question_answers")
$ws1.Range("B3").Value = "1"
$ws1.Range("B4").Value = "2"
$ws1.Range("B5").Value = "1"
$ws1.Range("B6").Value = "1"
$ws1.Range("B8").Value = "4"
$ws1.Range("B9").Value = "2"
$ws1.Range("B10").Value = "1"
$ws1.Range("B11").Value = "1"
$ws1.Range("B14").Value = "5"
$ws1.Range("B15").Value = "3"
$ws1.Range("B16").Value = "1"
$ws1.Range("B17").Value = "1"
$ws1.Range("B18").Value = "4"
$ws1.Range("B19").Value = "1"
$ws1.Range("B20").Value = "1"
$ws1.Range("B21").Value = "4"
$ws1.Range("B22").Value = "4"
$ws1.Range("B23").Value = "3"
$ws1.Range("B24").Value = "5"
$ws1.Range("B25").Value = "1"
$ws1.Range("B26").Value = "1"
$ws1.Range("B27").Value = "4"
$ws1.Range("B28").Value = "3"
$ws1.Range("B29").Value = "2"
$ws1.Range("B31").Value = "3"
$ws1.Range("B32").Value = "5"
$ws1.Range("B33").Value = "5"
$ws1.Range("B34").Value = "5"
$ws1.Range("B35").Value = "4"
$ws1.Range("B36").Value = "2"
$ws1.Range("B37").Value = "4"
$ws1.Range("B38").Value = "1"
$ws1.Range("B39").Value = "1"
$ws1.Range("B40").Value = "3"
$ws1.Range("B42").Value = "4"
$ws1.Range("B43").Value = "4"
$ws1.Range("B44").Value = "2"
$ws1.Range("B45").Value = "5"
$ws1.Range("B46").Value = "3"
$ws1.Range("B47").Value = "2"
$ws1.Range("B48").Value = "1"
$ws1.Range("B50").Value = "4"
$ws1.Range("B52").Value = "5"
$ws1.Range("B53").Value = "2"
$ws1.Range("B54").Value = "4"
$ws1.Range("B57").Value = "3"
$ws1.Range("B58").Value = "2"
$ws1.Range("B59").Value = "2"
$ws1.Range("B60").Value = "5"
$ws1.Range("B61").Value = "3"
$ws1.Range("B65").Value = "2"
$ws1.Range("B66").Value = "5"
$ws1.Range("B67").Value = "1"
$ws1.Range("B68").Value = "5"
$ws1.Range("B69").Value = "5"
$ws1.Range("B70").Value = "4"
$ws1.Range("B71").Value = "3"
$ws1.Range("B72").Value = "3"
$ws1.Range("B73").Value = "4"
$ws1.Range("B74").Value = "3"
$ws1.Range("B75").Value = "1"
$ws1.Range("B76").Value = "3"
$ws1.Range("B77").Value = "4"
$ws1.Range("B78").Value = "3"
$ws1.Range("B79").Value = "3"
$ws1.Range("B80").Value = "1"
$ws1.Range("B81").Value = "5"
$ws1.Range("B83").Value = "3"
$ws1.Range("B84").Value = "1"
$ws1.Range("B85").Value = "5"
$ws1.Range("B86").Value = "4"
$ws1.Range("B87").Value = "1"
$ws1.Range("B88").Value = "1"
$ws1.Range("B89").Value = "2"
$ws1.Range("B90").Value = "5"
$ws1.Range("B91").Value = "2"

# Sheet 2: outputs - update column B (numeric score values)
$ws2 = $wb.Worksheets.Item("outputs")
$ws2.Range("B2").Value = 271
$ws2.Range("B3").Value = 16
$ws2.Range("B5").Value = 21
$ws2.Range("B6").Value = 15
$ws2.Range("B7").Value = 18
$ws2.Range("B8").Value = 21
$ws2.Range("B9").Value = 16
$ws2.Range("B10").Value = 18
$ws2.Range("B11").Value = 21
$ws2.Range("B12").Value = 23
$ws2.Range("B13").Value = 19
$ws2.Range("B14").Value = 19
$ws2.Range("B15").Value = 18
$ws2.Range("B16").Value = 13
$ws2.Range("B17").Value = 15
$ws2.Range("B18").Value = 98
$ws2.Range("B19").Value = 49
$ws2.Range("B20").Value = 51
$ws2.Range("B21").Value = 39
$ws2.Range("B22").Value = 34
